# "Add files via upload" — populate the previously-empty Sheet1 with the
# small Name/Score/Notes roster table (C3:E7) and leave the selection
# resting one row below the data (D8), matching the uploaded workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("C3").Value = "Name "
$ws.Range("D3").Value = "Score"
$ws.Range("E3").Value = "Notes"

# Data rows (Notes column intentionally left blank, as in the source file)
$ws.Range("C4").Value = "Adam"
$ws.Range("D4").Value = 78

$ws.Range("C5").Value = "Bella"
$ws.Range("D5").Value = 98

$ws.Range("C6").Value = "Cindy"
$ws.Range("D6").Value = 60

$ws.Range("C7").Value = "David"
$ws.Range("D7").Value = 81

# Leave the cursor where the uploaded workbook had it: one row under the
# last data row, in the Score column.
[void]$ws.Range("D8").Select()

# Best-effort: restore the author's window size recorded in the workbook
# view (harmless no-op on hosts that don't wire this up to bookViews).
$win = $excel.ActiveWindow
$win.Width = 22260
$win.Height = 12650
